$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Stamp the "About" sheet with a date in cell C1 (serial 44307 = 2021-04-21),
# formatted using the built-in short-date number format.
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
